# Census data update: rename "White" column header to "White (%)" and
# backfill the previously-missing "Native Hawaiian (%)" values (column L)
# for the 4 states whose row was short one cell (ME, MI, VT, WV) - this
# had pushed their remaining columns (M/N/O) one slot to the left, leaving
# the White (%) figure absent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the last column header from "White" to "White (%)"
$ws.Range("O1").Value = "White (%)"

# Maine (row 21): insert missing Native Hawaiian (%) = 0, shift remaining values right
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = 1.7
$ws.Range("N21").Value = 1.6
$ws.Range("O21").Value = 93.5

# Michigan (row 24): same fix
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 2.4
$ws.Range("N24").Value = 5
$ws.Range("O24").Value = 75.4

# Vermont (row 47): same fix
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = 1.9
$ws.Range("N47").Value = 1.9
$ws.Range("O47").Value = 93.1

# West Virginia (row 50): same fix
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = 1.7
$ws.Range("N50").Value = 1.5
$ws.Range("O50").Value = 92.3

# Leave the selection where the author last clicked before saving
$ws.Range("S52").Select()
